$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue "D2" "36.678.92"
Set-TextValue "E2" "  -1.92%  "

Set-TextValue "D3" "2.018.59"
Set-TextValue "E3" "  +0.11%  "

Set-TextValue "D5" "235.18"
Set-TextValue "E5" "  -9.73%  "

Set-TextValue "D6" "0.601"
Set-TextValue "E6" "  -3.00%  "

Set-TextValue "D8" "54.97"
Set-TextValue "E8" "  -3.82%  "

Set-TextValue "D9" "0.371"
Set-TextValue "E9" "  -3.44%  "

Set-TextValue "D10" "57.61"
Set-TextValue "E10" "  +2.69%  "

Set-TextValue "D11" "0.0749"
Set-TextValue "E11" "  -3.69%  "

Set-TextValue "E12" "  -0.14%  "

Set-TextValue "D13" "2.304.64"
Set-TextValue "E13" "  -0.35%  "

Set-TextValue "D14" "14.18"
Set-TextValue "E14" "  -1.25%  "

Set-TextValue "D15" "20.23"
Set-TextValue "E15" "  -6.37%  "

Set-TextValue "D16" "0.759"
Set-TextValue "E16" "  -5.14%  "

Set-TextValue "E17" "  -2.81%  "

Set-TextValue "D18" "2.010.34"
Set-TextValue "E18" "  -1.13%  "

Set-TextValue "D19" "36.931.71"
Set-TextValue "E19" "  -1.10%  "

Set-TextValue "D20" "67.84"
Set-TextValue "E20" "  -3.25%  "

Set-TextValue "D21" "0.0₃0799"
Set-TextValue "E21" "  -4.83%  "

Set-TextValue "D22" "5.34"
Set-TextValue "E22" "  +3.97%  "

Set-TextValue "D23" "221.44"
Set-TextValue "E23" "  -5.30%  "

Set-TextValue "E24" "  +0.11%  "

Set-TextValue "E25" "  +2.41%  "

Set-TextValue "D26" "2.40"
Set-TextValue "E26" "  -7.98%  "

Set-TextValue "D27" "163.37"
Set-TextValue "E27" "  -0.96%  "

Set-TextValue "D28" "8.67"
Set-TextValue "E28" "  -3.65%  "

Set-TextValue "E29" "  +4.40%  "

Set-TextValue "E30" "  -1.85%  "

Set-TextValue "D31" "18.87"
Set-TextValue "E31" "  -3.83%  "

Set-TextValue "E32" "  -2.58%  "

Set-TextValue "D33" "4.39"
Set-TextValue "E33" "  -4.97%  "

Set-TextValue "D34" "0.0605"
Set-TextValue "E34" "  -5.95%  "

Set-TextValue "D35" "2.42"
Set-TextValue "E35" "  +1.66%  "

Set-TextValue "D36" "4.25"
Set-TextValue "E36" "  -6.11%  "

Set-TextValue "E37" "  +0.04%  "

Set-TextValue "E38" "  -2.30%  "

Set-TextValue "D39" "3.31"
Set-TextValue "E39" "  -1.61%  "

Set-TextValue "D40" "5.76"
Set-TextValue "E40" "  +4.47%  "

Set-TextValue "D41" "2.99"
Set-TextValue "E41" "  -2.24%  "

Set-TextValue "D42" "1.461.70"
Set-TextValue "E42" "  +1.79%  "

Set-TextValue "D43" "0.0931"
Set-TextValue "E43" "  -0.08%  "

Set-TextValue "D44" "0.0205"
Set-TextValue "E44" "  -3.95%  "

Set-TextValue "B45" "TrustWalletToken"
Set-TextValue "C45" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D45" "1.11"
Set-TextValue "E45" "  -7.69%  "

Set-TextValue "B46" "Aave"
Set-TextValue "C46" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D46" "90.09"
Set-TextValue "E46" "  +0.20%  "

Set-TextValue "D47" "15.36"
Set-TextValue "E47" "  -2.86%  "

Set-TextValue "D48" "1.00"
Set-TextValue "E48" "  -3.01%  "

Set-TextValue "D49" "3.85"
Set-TextValue "E49" "  +26.40%  "

Set-TextValue "D50" "2.89"
Set-TextValue "E50" "  -1.45%  "

Set-TextValue "D51" "6.85"
Set-TextValue "E51" "  -2.69%  "
